$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force plain-number-looking price strings to stay as text (avoid Excel
# auto-converting values like "1.003" into the number 1.003).
$textCells = @("D4","D5","D6","D7","D8","D9","D10","D11","D12","D14","D15","D16","D18","D19","D21","D22","D24","D25","D26","D27","D29","D30","D31","D33","D34","D35","D37","D39","D40","D41","D42","D43","D44","D45","D46","D47","D48","D49","D50","D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply updated coin data scraped this run
# Row 2
$ws.Range("D2").Value = '28.022.83'
$ws.Range("E2").Value = '  +2.16%  '

# Row 3
$ws.Range("D3").Value = '1.910.93'
$ws.Range("E3").Value = '  +2.49%  '

# Row 4
$ws.Range("D4").Value = '1.003'
$ws.Range("E4").Value = '  -0.93%  '

# Row 5
$ws.Range("D5").Value = '315.42'
$ws.Range("E5").Value = '  +1.27%  '

# Row 6
$ws.Range("D6").Value = '1.004'
$ws.Range("E6").Value = '  -0.79%  '

# Row 7
$ws.Range("D7").Value = '0.4806'
$ws.Range("E7").Value = '  +0.59%  '

# Row 8
$ws.Range("D8").Value = '0.3813'
$ws.Range("E8").Value = '  +1.41%  '

# Row 9
$ws.Range("D9").Value = '0.07367'
$ws.Range("E9").Value = '  +0.64%  '

# Row 10
$ws.Range("D10").Value = '0.9346'
$ws.Range("E10").Value = '  -0.04%  '

# Row 11
$ws.Range("D11").Value = '20.83'
$ws.Range("E11").Value = '  +0.83%  '

# Row 12
$ws.Range("D12").Value = '0.07777'
$ws.Range("E12").Value = '  -0.54%  '

# Row 13
$ws.Range("D13").Value = '1.881.52'
$ws.Range("E13").Value = '  +0.39%  '

# Row 14
$ws.Range("D14").Value = '5.501'
$ws.Range("E14").Value = '  +1.26%  '

# Row 15
$ws.Range("D15").Value = '6.642'
$ws.Range("E15").Value = '  +1.33%  '

# Row 16
$ws.Range("D16").Value = '91.99'
$ws.Range("E16").Value = '  +1.64%  '

# Row 17
$ws.Range("E17").Value = '  -0.88%  '

# Row 18
$ws.Range("D18").Value = '0.000008854'
$ws.Range("E18").Value = '  -0.22%  '

# Row 19
$ws.Range("D19").Value = '1.003'
$ws.Range("E19").Value = '  -0.78%  '

# Row 20
$ws.Range("D20").Value = '28.048.93'
$ws.Range("E20").Value = '  +2.04%  '

# Row 21
$ws.Range("D21").Value = '14.80'
$ws.Range("E21").Value = '  +0.45%  '

# Row 22
$ws.Range("D22").Value = '5.173'
$ws.Range("E22").Value = '  +1.18%  '

# Row 23
$ws.Range("D23").Value = '2.143.28'
$ws.Range("E23").Value = '  +1.04%  '

# Row 24
$ws.Range("D24").Value = '10.90'
$ws.Range("E24").Value = '  +1.95%  '

# Row 25
$ws.Range("D25").Value = '155.72'
$ws.Range("E25").Value = '  +0.23%  '

# Row 26
$ws.Range("D26").Value = '1.916'
$ws.Range("E26").Value = '  -1.16%  '

# Row 27
$ws.Range("D27").Value = '18.50'
$ws.Range("E27").Value = '  +0.19%  '

# Row 28
$ws.Range("E28").Value = '  +6.15%  '

# Row 29
$ws.Range("D29").Value = '117.00'
$ws.Range("E29").Value = '  +1.32%  '

# Row 30
$ws.Range("D30").Value = '4.971'
$ws.Range("E30").Value = '  +0.14%  '

# Row 31
$ws.Range("D31").Value = '0.08949'
$ws.Range("E31").Value = '  +0.60%  '

# Row 32
$ws.Range("E32").Value = '  -0.99%  '

# Row 33
$ws.Range("D33").Value = '1.267'
$ws.Range("E33").Value = '  +4.45%  '

# Row 34
$ws.Range("D34").Value = '0.7797'
$ws.Range("E34").Value = '  +3.44%  '

# Row 35
$ws.Range("D35").Value = '4.688'
$ws.Range("E35").Value = '  +1.90%  '

# Row 36
$ws.Range("E36").Value = '  -3.61%  '

# Row 37
$ws.Range("D37").Value = '0.02057'
$ws.Range("E37").Value = '  +1.29%  '

# Row 38
$ws.Range("E38").Value = '  -0.83%  '

# Row 39
$ws.Range("D39").Value = '0.05325'
$ws.Range("E39").Value = '  +1.26%  '

# Row 40
$ws.Range("B40").Value = 'TheSandbox'
$ws.Range("C40").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D40").Value = '0.5500'
$ws.Range("E40").Value = '  +3.49%  '

# Row 41
$ws.Range("B41").Value = 'MXToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D41").Value = '3.004'
$ws.Range("E41").Value = '  +0.50%  '

# Row 42
$ws.Range("D42").Value = '7.030'
$ws.Range("E42").Value = '  -0.61%  '

# Row 43
$ws.Range("D43").Value = '0.1530'
$ws.Range("E43").Value = '  +0.45%  '

# Row 44
$ws.Range("D44").Value = '8.498'
$ws.Range("E44").Value = '  -0.64%  '

# Row 45
$ws.Range("D45").Value = '10.65'
$ws.Range("E45").Value = '  +0.38%  '

# Row 46
$ws.Range("D46").Value = '0.4834'
$ws.Range("E46").Value = '  +0.69%  '

# Row 47
$ws.Range("D47").Value = '107.99'
$ws.Range("E47").Value = '  +5.05%  '

# Row 48
$ws.Range("D48").Value = '1.004'
$ws.Range("E48").Value = '  -0.83%  '

# Row 49
$ws.Range("D49").Value = '1.652'
$ws.Range("E49").Value = '  -0.26%  '

# Row 50
$ws.Range("D50").Value = '68.06'
$ws.Range("E50").Value = '  +1.19%  '

# Row 51
$ws.Range("D51").Value = '0.06071'
$ws.Range("E51").Value = '  -0.14%  '
